# Apply ETL-consolidator refactor data changes to the absenteeism sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 90192
$ws.Range("B2").Value = "Dra. Marcela Ramos"
$ws.Range("C2").Value = "Financeiro"
$ws.Range("E2").Value = 6
$ws.Range("F2").Value = 45099
$ws.Range("G2").Value = 9786.15

# Row 3
$ws.Range("A3").Value = 95112
$ws.Range("B3").Value = "Elisa Teixeira"
$ws.Range("C3").Value = "Operações"
$ws.Range("D3").Value = "Doença"
$ws.Range("E3").Value = 7
$ws.Range("F3").Value = 45102
$ws.Range("G3").Value = 8248.549999999999

# Row 4
$ws.Range("A4").Value = 43938
$ws.Range("B4").Value = "Dr. Vinicius Barros"
$ws.Range("C4").Value = "Jurídico"
$ws.Range("D4").Value = "Doença"
$ws.Range("E4").Value = 7
$ws.Range("F4").Value = 45098
$ws.Range("G4").Value = 4903.32

# Row 5
$ws.Range("A5").Value = 1203
$ws.Range("B5").Value = "Rafaela Gonçalves"
$ws.Range("C5").Value = "Operações"
$ws.Range("D5").Value = "Viagem de negócios"
$ws.Range("E5").Value = 4
$ws.Range("F5").Value = 45090
$ws.Range("G5").Value = 5928.88

# Row 6
$ws.Range("A6").Value = 27547
$ws.Range("B6").Value = "Bárbara Nogueira"
$ws.Range("C6").Value = "Financeiro"
$ws.Range("D6").Value = "Consulta médica"
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 45104
$ws.Range("G6").Value = 5166.6

# Row 7
$ws.Range("A7").Value = 66548
$ws.Range("B7").Value = "Dr. Emanuel Duarte"
$ws.Range("C7").Value = "Operações"
$ws.Range("D7").Value = "Consulta médica"
$ws.Range("E7").Value = 7
$ws.Range("F7").Value = 45099
$ws.Range("G7").Value = 10857.71

# Row 8
$ws.Range("A8").Value = 88219
$ws.Range("B8").Value = "Caio Jesus"
$ws.Range("C8").Value = "Recursos Humanos"
$ws.Range("D8").Value = "Doença"
$ws.Range("E8").Value = 7
$ws.Range("F8").Value = 45101
$ws.Range("G8").Value = 11509.21

# Row 9
$ws.Range("A9").Value = 84171
$ws.Range("B9").Value = "Luiz Felipe Peixoto"
$ws.Range("E9").Value = 5
$ws.Range("F9").Value = 45100
$ws.Range("G9").Value = 9551.07

# Row 10
$ws.Range("A10").Value = 92708
$ws.Range("B10").Value = "Ana Julia Mendes"
$ws.Range("C10").Value = "Atendimento ao Cliente"
$ws.Range("D10").Value = "Viagem de negócios"
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 45081
$ws.Range("G10").Value = 5976.08

# Row 11
$ws.Range("A11").Value = 40057
$ws.Range("B11").Value = "Marina da Rocha"
$ws.Range("C11").Value = "Jurídico"
$ws.Range("E11").Value = 5
$ws.Range("F11").Value = 45084
$ws.Range("G11").Value = 8156.22
